# feat: add 2022-Q3 data
#
# Inserts a new "2022-Q3" worksheet (right after "总计") populated with the
# fund-holding detail table, and updates the "总计" (summary) worksheet so its
# small history table gains a new first data row for 2022-Q3 (shifting the
# existing rows down) plus the 2020-Q4 row that falls off the bottom.

$wb    = $excel.ActiveWorkbook
$total = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1. Insert the new "2022-Q3" sheet right after "总计", before the old
#    "2021-Q3" sheet (mirrors the diff: sheetId 2 becomes the new sheet, all
#    the old quarters shift one slot to the right/down).
# ---------------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add($null, $total)
$newSheet.Name = "2022-Q3"

# Use the (now shifted) "2021-Q3" sheet purely as a formatting template so the
# new sheet's header row / index column match the existing look (bold,
# centered, bordered). Column A is left alone on row 1 (no sheet has an A1
# cell), so only B1:H1 is copied for the header.
$tmpl = $wb.Worksheets.Item("2021-Q3")
$tmpl.Range("B1:H1").Copy($newSheet.Range("B1"))
$tmpl.Range("A2").Copy($newSheet.Range("A2:A19"))

# Header text (D1 differs from the other quarter sheets: "基金规模" not "基金金额").
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Column A is just the 0-based row index (0..17).
for ($i = 0; $i -le 17; $i++) {
    $newSheet.Cells.Item($i + 2, 1).Value = $i
}

# ---------------------------------------------------------------------------
# 2. Fund-holding detail rows for 2022-Q3.
# ---------------------------------------------------------------------------
$rows = @(
    @("009011", "华夏睿阳一年持有期混合",           "15.05", "75.73", "1.87", "0.2814", 10),
    @("001195", "工银农业产业股票",                 "5.53",  "80.67", "4.24", "0.2345", 3),
    @("161609", "融通动力先锋混合",                 "6.67",  "86.12", "3.07", "0.2048", 6),
    @("010994", "博时创新经济混合A",                "3.35",  "91.74", "5.44", "0.1822", 4),
    @("004683", "建信高端医疗股票A",                "2.88",  "85.77", "3.37", "0.0971", 6),
    @("001152", "融通新区域新经济灵活配置混合",      "2.64",  "86.64", "3.05", "0.0805", 7),
    @("004050", "华夏新锦升灵活配置混合A",           "1.22",  "64.62", "5.82", "0.0710", 1),
    @("014781", "建信兴衡优选一年持有期混合A",       "1.77",  "46.74", "3.52", "0.0623", 2),
    @("015056", "百嘉百盛混合",                     "1.21",  "62.90", "3.26", "0.0394", 2),
    @("014782", "建信兴衡优选一年持有期混合C",       "0.70",  "46.74", "3.52", "0.0246", 2),
    @("010995", "博时创新经济混合C",                "0.38",  "91.74", "5.44", "0.0207", 4),
    @("006225", "人保量化基本面混合A",               "0.47",  "90.98", "3.86", "0.0181", 5),
    @("004051", "华夏新锦升灵活配置混合C",           "0.06",  "64.62", "5.82", "0.0035", 1),
    @("001657", "长安鑫富领先灵活配置混合",           "0.06",  "49.60", "3.22", "0.0019", 5),
    @("006226", "人保量化基本面混合C",               "0.04",  "90.98", "3.86", "0.0015", 5),
    @("166107", "信澳量化多因子混合（LOF）A",         "0.05",  "28.39", "0.38", "0.0002", 9),
    @("166108", "信澳量化多因子混合（LOF）C",         "0.06",  "28.39", "0.38", "0.0002", 9),
    @("016352", "建信高端医疗股票C",                "0.00",  "85.77", "3.37", $null,     6)
)

$r = 2
foreach ($row in $rows) {
    $code  = $row[0]
    $name  = $row[1]
    $size  = $row[2]
    $pos   = $row[3]
    $pct   = $row[4]
    $mktv  = $row[5]
    $rank  = $row[6]

    $codeCell = $newSheet.Cells.Item($r, 2)
    $codeCell.NumberFormat = "@"
    $codeCell.Value = $code

    $newSheet.Cells.Item($r, 3).Value = $name

    $sizeCell = $newSheet.Cells.Item($r, 4)
    $sizeCell.NumberFormat = "@"
    $sizeCell.Value = $size

    $posCell = $newSheet.Cells.Item($r, 5)
    $posCell.NumberFormat = "@"
    $posCell.Value = $pos

    $pctCell = $newSheet.Cells.Item($r, 6)
    $pctCell.NumberFormat = "@"
    $pctCell.Value = $pct

    if ($null -eq $mktv) {
        $newSheet.Cells.Item($r, 7).Value = 0
    } else {
        $mktvCell = $newSheet.Cells.Item($r, 7)
        $mktvCell.NumberFormat = "@"
        $mktvCell.Value = $mktv
    }

    $newSheet.Cells.Item($r, 8).Value = $rank

    $r++
}

# ---------------------------------------------------------------------------
# 3. Update the "总计" (summary) sheet: insert a new row 2 for 2022-Q3 and
#    append a new trailing row for 2020-Q4 (the quarter data shifts down one
#    row; column A is a plain 0-based row index so it is restored to
#    0,1,2,3,4 afterwards rather than shifting along with B:D).
# ---------------------------------------------------------------------------
$total.Rows.Item(2).Insert()
$total.Range("A2:D2").Style = "Normal"
$total.Range("A3").Copy($total.Range("A2"))
$total.Range("B2").Value = "2022-Q3"
$total.Range("C2").Value = 18
$total.Range("D2").Value = 1.32

$total.Range("A5").Copy($total.Range("A6"))
$total.Range("B6").Value = "2020-Q4"
$total.Range("C6").Value = 6
$total.Range("D6").Value = 0.65

# Column A is untouched by the diff for the pre-existing rows (2-5): restore
# the original 0,1,2,3 index values instead of letting them shift along with
# the row insert, then set 4 for the brand-new trailing row.
$total.Range("A2").Value = 0
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3
$total.Range("A6").Value = 4

# ---------------------------------------------------------------------------
# 4. Restore "2020-Q4" as the selected/active sheet, matching the original
#    workbook (unrelated to the new data, but otherwise untouched by the diff).
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("2020-Q4").Activate()
